$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cryptocurrency table values are plain text in the workbook (they were
# written as inline/shared strings, e.g. "26.191.22", "0.2510", "1.015").
# Some of those look like plain decimal numbers to Excel, so for those cells we
# force the Text number format ('@') before writing the value, which prevents
# Excel from silently converting the string into a number and stripping
# significant leading/trailing zeros. Cells whose text would never be
# reinterpreted (coin names, URLs, multi-dot price strings, padded percentages)
# are written directly, matching the original file untouched.

$ws.Range('D2').Value = '26.191.22'
$ws.Range('E2').Value = '  -1.05%  '

$ws.Range('D3').Value = '1.826.48'
$ws.Range('E3').Value = '  -0.76%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.015'
$ws.Range('E4').Value = '  +1.38%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.66'
$ws.Range('E5').Value = '  -8.03%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.58%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5098'
$ws.Range('E7').Value = '  -2.91%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2510'
$ws.Range('E8').Value = '  -21.44%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06295'
$ws.Range('E9').Value = '  -7.26%  '

$ws.Range('D10').Value = '1.832.22'
$ws.Range('E10').Value = '  -0.08%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06847'
$ws.Range('E11').Value = '  -11.60%  '

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.52'
$ws.Range('E12').Value = '  -22.65%  '

$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '80.17'
$ws.Range('E13').Value = '  -8.59%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.373'
$ws.Range('E14').Value = '  -12.79%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5597'
$ws.Range('E15').Value = '  -28.81%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.022'
$ws.Range('E16').Value = '  +2.11%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.36%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '26.238.86'
$ws.Range('E18').Value = '  -0.97%  '

$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.56'
$ws.Range('E19').Value = '  -16.54%  '

$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.071.72'
$ws.Range('E20').Value = '  -0.07%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006221'
$ws.Range('E21').Value = '  -21.72%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.972'
$ws.Range('E22').Value = '  -14.14%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.272'
$ws.Range('E23').Value = '  -11.80%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.933'
$ws.Range('E24').Value = '  -15.38%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '135.01'
$ws.Range('E25').Value = '  -4.39%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.607'
$ws.Range('E26').Value = '  -4.30%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.863'
$ws.Range('E27').Value = '  -14.40%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '14.84'
$ws.Range('E28').Value = '  -12.39%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '99.71'
$ws.Range('E29').Value = '  -10.76%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08520'
$ws.Range('E30').Value = '  -1.90%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.693'
$ws.Range('E31').Value = '  -11.17%  '

$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.848'
$ws.Range('E32').Value = '  -0.33%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04366'
$ws.Range('E33').Value = '  -10.33%  '

$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.208'
$ws.Range('E34').Value = '  -21.24%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.049'
$ws.Range('E35').Value = '  -7.57%  '

$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.031'
$ws.Range('E36').Value = '  -1.94%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6245'
$ws.Range('E37').Value = '  -14.19%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.032'
$ws.Range('E38').Value = '  -9.52%  '

$ws.Range('B39').Value = 'PaxDollar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.006'
$ws.Range('E39').Value = '  +0.51%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01501'
$ws.Range('E40').Value = '  -14.44%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8240'
$ws.Range('E41').Value = '  -7.60%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.569'
$ws.Range('E42').Value = '  -6.05%  '

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.66'
$ws.Range('E43').Value = '  -8.98%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3872'
$ws.Range('E44').Value = '  -18.95%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.495'
$ws.Range('E45').Value = '  -5.08%  '

$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05386'
$ws.Range('E46').Value = '  -7.90%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.305'
$ws.Range('E47').Value = '  -17.93%  '

$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '31.16'
$ws.Range('E48').Value = '  -10.62%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.03'
$ws.Range('E49').Value = '  -11.00%  '

$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.013'
$ws.Range('E50').Value = '  +0.81%  '

$ws.Range('B51').Value = 'TrueUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('E51').Value = '  +0.17%  '
